$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers to Excel's
# auto-detection need to be forced to Text format first, so they are
# written as text (matching the source data) instead of being parsed
# into floating point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "54.385.13"
$ws.Range("E2").Value = "  +5.24%  "
$ws.Range("D3").Value = "3.176.67"
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "398.03"
$ws.Range("E5").Value = "  +2.58%  "
$ws.Range("D6").Value = "109.36"
$ws.Range("E6").Value = "  +5.81%  "
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").Value = "  +4.76%  "
$ws.Range("D10").Value = "38.90"
$ws.Range("E10").Value = "  +4.96%  "
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").Value = "0.0884"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("D13").Value = "3.673.30"
$ws.Range("E13").Value = "  +2.98%  "
$ws.Range("D14").Value = "19.18"
$ws.Range("E14").Value = "  +2.33%  "
$ws.Range("D15").Value = "8.08"
$ws.Range("E15").Value = "  +3.15%  "
$ws.Range("E16").Value = "  +8.29%  "
$ws.Range("D17").Value = "3.175.93"
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("D18").Value = "10.51"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").Value = "54.330.14"
$ws.Range("E19").Value = "  +4.95%  "
$ws.Range("D20").Value = "3.30"
$ws.Range("E20").Value = "  +4.39%  "
$ws.Range("D21").Value = "12.89"
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("D22").Value = "0.0000100"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").Value = "71.33"
$ws.Range("E23").Value = "  +1.38%  "
$ws.Range("D24").Value = "272.58"
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").Value = "3.24"
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("D26").Value = "8.03"
$ws.Range("E26").Value = "  -2.40%  "
$ws.Range("D27").Value = "27.76"
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").Value = "0.171"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  +3.90%  "
$ws.Range("D32").Value = "11.03"
$ws.Range("E32").Value = "  +7.04%  "
$ws.Range("D33").Value = "0.0499"
$ws.Range("E33").Value = "  +10.31%  "
$ws.Range("D34").Value = "36.95"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").Value = "50.57"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").Value = "3.63"
$ws.Range("E37").Value = "  +8.77%  "
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "2.86"
$ws.Range("E39").Value = "  +11.34%  "
$ws.Range("D40").Value = "4.13"
$ws.Range("E40").Value = "  +10.47%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "17.43"
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "0.292"
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("D44").Value = "129.61"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("E45").Value = "  +1.05%  "
$ws.Range("D46").Value = "22.33"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("D47").Value = "2.44"
$ws.Range("E47").Value = "  -0.74%  "
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D49").Value = "2.086.96"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("E50").Value = "  +7.05%  "
$ws.Range("D51").Value = "0.0500"
$ws.Range("E51").Value = "  +11.36%  "
